$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 326
$ws.Cells.Item(326, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(326, 2).Value = 8
$ws.Cells.Item(326, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(326, 4).Value = "D. RODRIGUEZ GARCIA"
$ws.Cells.Item(326, 5).Value = "C. D. MENSAJERO ISLA DE LA PALMA"
$ws.Cells.Item(326, 6).Value = 2.08
$ws.Cells.Item(326, 7).Value = 124.8
$ws.Cells.Item(326, 8).Value = 0
$ws.Cells.Item(326, 9).Value = 1
$ws.Cells.Item(326, 10).Value = 0
$ws.Cells.Item(326, 11).Value = 1
$ws.Cells.Item(326, 12).Value = 0
$ws.Cells.Item(326, 13).Value = 0
$ws.Cells.Item(326, 14).Value = 0
$ws.Cells.Item(326, 15).Value = 0
$ws.Cells.Item(326, 16).Value = 0
$ws.Cells.Item(326, 17).Value = 0
$ws.Cells.Item(326, 18).Value = 0
$ws.Cells.Item(326, 19).Value = 0
$ws.Cells.Item(326, 20).Value = 0
$ws.Cells.Item(326, 21).Value = 0
$ws.Cells.Item(326, 22).Value = 0
$ws.Cells.Item(326, 23).Value = 21.01
$ws.Cells.Item(326, 24).Value = -3
$ws.Cells.Item(326, 25).Value = -55.97

# Row 327
$ws.Cells.Item(327, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(327, 2).Value = 8
$ws.Cells.Item(327, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(327, 4).Value = "P. RODRIGUEZ RIVERO"
$ws.Cells.Item(327, 5).Value = "C. D. MENSAJERO ISLA DE LA PALMA"
$ws.Cells.Item(327, 6).Value = 2.08
$ws.Cells.Item(327, 7).Value = 124.8
$ws.Cells.Item(327, 8).Value = 0
$ws.Cells.Item(327, 9).Value = 1
$ws.Cells.Item(327, 10).Value = 0
$ws.Cells.Item(327, 11).Value = 1
$ws.Cells.Item(327, 12).Value = 0
$ws.Cells.Item(327, 13).Value = 2
$ws.Cells.Item(327, 14).Value = 2
$ws.Cells.Item(327, 15).Value = 0
$ws.Cells.Item(327, 16).Value = 0
$ws.Cells.Item(327, 17).Value = 0
$ws.Cells.Item(327, 18).Value = 0
$ws.Cells.Item(327, 19).Value = 0
$ws.Cells.Item(327, 20).Value = 0
$ws.Cells.Item(327, 21).Value = 0
$ws.Cells.Item(327, 22).Value = 0
$ws.Cells.Item(327, 23).Value = 39.5
$ws.Cells.Item(327, 24).Value = -3
$ws.Cells.Item(327, 25).Value = -55.97

# Row 328
$ws.Cells.Item(328, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(328, 2).Value = 8
$ws.Cells.Item(328, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(328, 4).Value = "J. RIES"
$ws.Cells.Item(328, 5).Value = "C. D. MENSAJERO ISLA DE LA PALMA"
$ws.Cells.Item(328, 6).Value = 2.08
$ws.Cells.Item(328, 7).Value = 124.8
$ws.Cells.Item(328, 8).Value = 0
$ws.Cells.Item(328, 9).Value = 0
$ws.Cells.Item(328, 10).Value = 0
$ws.Cells.Item(328, 11).Value = 0
$ws.Cells.Item(328, 12).Value = 0
$ws.Cells.Item(328, 13).Value = 0
$ws.Cells.Item(328, 14).Value = 0
$ws.Cells.Item(328, 17).Value = 0
$ws.Cells.Item(328, 18).Value = 1
$ws.Cells.Item(328, 19).Value = 0
$ws.Cells.Item(328, 20).Value = 1
$ws.Cells.Item(328, 21).Value = 0
$ws.Cells.Item(328, 22).Value = 1
$ws.Cells.Item(328, 23).Value = 21.01
$ws.Cells.Item(328, 24).Value = -3
$ws.Cells.Item(328, 25).Value = -55.97

# Row 329
$ws.Cells.Item(329, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(329, 2).Value = 8
$ws.Cells.Item(329, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(329, 4).Value = "O. PEÑA LOPEZ"
$ws.Cells.Item(329, 5).Value = "C. D. MENSAJERO ISLA DE LA PALMA"
$ws.Cells.Item(329, 6).Value = 2.08
$ws.Cells.Item(329, 7).Value = 124.8
$ws.Cells.Item(329, 8).Value = 0
$ws.Cells.Item(329, 9).Value = 2
$ws.Cells.Item(329, 10).Value = 1
$ws.Cells.Item(329, 11).Value = 0
$ws.Cells.Item(329, 12).Value = 0
$ws.Cells.Item(329, 13).Value = 2
$ws.Cells.Item(329, 14).Value = 0
$ws.Cells.Item(329, 15).Value = 0.5
$ws.Cells.Item(329, 16).Value = 0
$ws.Cells.Item(329, 17).Value = 0
$ws.Cells.Item(329, 18).Value = 0
$ws.Cells.Item(329, 19).Value = 1
$ws.Cells.Item(329, 20).Value = 2
$ws.Cells.Item(329, 21).Value = 2
$ws.Cells.Item(329, 22).Value = 0
$ws.Cells.Item(329, 23).Value = 60.5
$ws.Cells.Item(329, 24).Value = -3
$ws.Cells.Item(329, 25).Value = -55.97

# Row 330
$ws.Cells.Item(330, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(330, 2).Value = 8
$ws.Cells.Item(330, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(330, 4).Value = "M. NIANG"
$ws.Cells.Item(330, 5).Value = "C. D. MENSAJERO ISLA DE LA PALMA"
$ws.Cells.Item(330, 6).Value = 1.63
$ws.Cells.Item(330, 7).Value = 97.8
$ws.Cells.Item(330, 8).Value = 0
$ws.Cells.Item(330, 9).Value = 0
$ws.Cells.Item(330, 10).Value = 0
$ws.Cells.Item(330, 11).Value = 0
$ws.Cells.Item(330, 12).Value = 0
$ws.Cells.Item(330, 13).Value = 0
$ws.Cells.Item(330, 14).Value = 0
$ws.Cells.Item(330, 17).Value = 0
$ws.Cells.Item(330, 18).Value = 0
$ws.Cells.Item(330, 19).Value = 0
$ws.Cells.Item(330, 20).Value = 0
$ws.Cells.Item(330, 21).Value = 0
$ws.Cells.Item(330, 22).Value = 0
$ws.Cells.Item(330, 23).Value = 0
$ws.Cells.Item(330, 24).Value = -2
$ws.Cells.Item(330, 25).Value = -48.45

# Row 331
$ws.Cells.Item(331, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(331, 2).Value = 8
$ws.Cells.Item(331, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(331, 4).Value = "A. APARICIO IZQUIERDO"
$ws.Cells.Item(331, 5).Value = "C. D. MENSAJERO ISLA DE LA PALMA"
$ws.Cells.Item(331, 6).Value = 0.45
$ws.Cells.Item(331, 7).Value = 27
$ws.Cells.Item(331, 8).Value = 0
$ws.Cells.Item(331, 9).Value = 0
$ws.Cells.Item(331, 10).Value = 0
$ws.Cells.Item(331, 11).Value = 0
$ws.Cells.Item(331, 12).Value = 0
$ws.Cells.Item(331, 13).Value = 0
$ws.Cells.Item(331, 14).Value = 0
$ws.Cells.Item(331, 17).Value = 0
$ws.Cells.Item(331, 18).Value = 0
$ws.Cells.Item(331, 19).Value = 0
$ws.Cells.Item(331, 20).Value = 0
$ws.Cells.Item(331, 21).Value = 0
$ws.Cells.Item(331, 22).Value = 0
$ws.Cells.Item(331, 23).Value = 0
$ws.Cells.Item(331, 24).Value = -1
$ws.Cells.Item(331, 25).Value = -72.73

# Row 332
$ws.Cells.Item(332, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(332, 2).Value = 8
$ws.Cells.Item(332, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(332, 4).Value = "G. DIAZ MONTERO"
$ws.Cells.Item(332, 5).Value = "C.B. TRES CANTOS"
$ws.Cells.Item(332, 6).Value = 2.08
$ws.Cells.Item(332, 7).Value = 124.8
$ws.Cells.Item(332, 8).Value = 0
$ws.Cells.Item(332, 9).Value = 2
$ws.Cells.Item(332, 10).Value = 0
$ws.Cells.Item(332, 11).Value = 0
$ws.Cells.Item(332, 12).Value = 0
$ws.Cells.Item(332, 13).Value = 0
$ws.Cells.Item(332, 14).Value = 0
$ws.Cells.Item(332, 15).Value = 0
$ws.Cells.Item(332, 16).Value = 0
$ws.Cells.Item(332, 17).Value = 0
$ws.Cells.Item(332, 18).Value = 1
$ws.Cells.Item(332, 19).Value = 0
$ws.Cells.Item(332, 20).Value = 2
$ws.Cells.Item(332, 21).Value = 1
$ws.Cells.Item(332, 22).Value = 1
$ws.Cells.Item(332, 23).Value = 60
$ws.Cells.Item(332, 24).Value = 3
$ws.Cells.Item(332, 25).Value = 55.97

# Row 333
$ws.Cells.Item(333, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(333, 2).Value = 8
$ws.Cells.Item(333, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(333, 4).Value = "N. MAIGA"
$ws.Cells.Item(333, 5).Value = "C.B. TRES CANTOS"
$ws.Cells.Item(333, 6).Value = 2.08
$ws.Cells.Item(333, 7).Value = 124.8
$ws.Cells.Item(333, 8).Value = 0
$ws.Cells.Item(333, 9).Value = 1
$ws.Cells.Item(333, 10).Value = 1
$ws.Cells.Item(333, 11).Value = 0
$ws.Cells.Item(333, 12).Value = 0
$ws.Cells.Item(333, 13).Value = 0
$ws.Cells.Item(333, 14).Value = 0
$ws.Cells.Item(333, 15).Value = 1
$ws.Cells.Item(333, 16).Value = 0
$ws.Cells.Item(333, 17).Value = 0
$ws.Cells.Item(333, 18).Value = 1
$ws.Cells.Item(333, 19).Value = 0
$ws.Cells.Item(333, 20).Value = 1
$ws.Cells.Item(333, 21).Value = 1
$ws.Cells.Item(333, 22).Value = 0
$ws.Cells.Item(333, 23).Value = 40
$ws.Cells.Item(333, 24).Value = 3
$ws.Cells.Item(333, 25).Value = 55.97

# Row 334
$ws.Cells.Item(334, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(334, 2).Value = 8
$ws.Cells.Item(334, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(334, 4).Value = "J. DOMINGUEZ LARRE"
$ws.Cells.Item(334, 5).Value = "C.B. TRES CANTOS"
$ws.Cells.Item(334, 6).Value = 2.08
$ws.Cells.Item(334, 7).Value = 124.8
$ws.Cells.Item(334, 8).Value = 0
$ws.Cells.Item(334, 9).Value = 1
$ws.Cells.Item(334, 10).Value = 1
$ws.Cells.Item(334, 11).Value = 1
$ws.Cells.Item(334, 12).Value = 1
$ws.Cells.Item(334, 13).Value = 0
$ws.Cells.Item(334, 14).Value = 0
$ws.Cells.Item(334, 15).Value = 1.5
$ws.Cells.Item(334, 16).Value = 0
$ws.Cells.Item(334, 17).Value = 0
$ws.Cells.Item(334, 18).Value = 0
$ws.Cells.Item(334, 19).Value = 0
$ws.Cells.Item(334, 20).Value = 0
$ws.Cells.Item(334, 21).Value = 0
$ws.Cells.Item(334, 22).Value = 0
$ws.Cells.Item(334, 23).Value = 20
$ws.Cells.Item(334, 24).Value = 3
$ws.Cells.Item(334, 25).Value = 55.97

# Row 335
$ws.Cells.Item(335, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(335, 2).Value = 8
$ws.Cells.Item(335, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(335, 4).Value = "D. GONZALEZ LONGARELA"
$ws.Cells.Item(335, 5).Value = "C.B. TRES CANTOS"
$ws.Cells.Item(335, 6).Value = 1.63
$ws.Cells.Item(335, 7).Value = 97.8
$ws.Cells.Item(335, 8).Value = 0
$ws.Cells.Item(335, 9).Value = 0
$ws.Cells.Item(335, 10).Value = 0
$ws.Cells.Item(335, 11).Value = 0
$ws.Cells.Item(335, 12).Value = 0
$ws.Cells.Item(335, 13).Value = 0
$ws.Cells.Item(335, 14).Value = 0
$ws.Cells.Item(335, 17).Value = 0
$ws.Cells.Item(335, 18).Value = 0
$ws.Cells.Item(335, 19).Value = 0
$ws.Cells.Item(335, 20).Value = 0
$ws.Cells.Item(335, 21).Value = 0
$ws.Cells.Item(335, 22).Value = 0
$ws.Cells.Item(335, 23).Value = 0
$ws.Cells.Item(335, 24).Value = 2
$ws.Cells.Item(335, 25).Value = 48.45

# Row 336
$ws.Cells.Item(336, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(336, 2).Value = 8
$ws.Cells.Item(336, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(336, 4).Value = "J. ATIENZA PEREA"
$ws.Cells.Item(336, 5).Value = "C.B. TRES CANTOS"
$ws.Cells.Item(336, 6).Value = 1.08
$ws.Cells.Item(336, 7).Value = 64.80000000000001
$ws.Cells.Item(336, 8).Value = 0
$ws.Cells.Item(336, 9).Value = 1
$ws.Cells.Item(336, 10).Value = 1
$ws.Cells.Item(336, 11).Value = 0
$ws.Cells.Item(336, 12).Value = 0
$ws.Cells.Item(336, 13).Value = 0
$ws.Cells.Item(336, 14).Value = 0
$ws.Cells.Item(336, 15).Value = 1
$ws.Cells.Item(336, 16).Value = 0
$ws.Cells.Item(336, 17).Value = 0
$ws.Cells.Item(336, 18).Value = 0
$ws.Cells.Item(336, 19).Value = 1
$ws.Cells.Item(336, 20).Value = 0
$ws.Cells.Item(336, 21).Value = 0
$ws.Cells.Item(336, 22).Value = 0
$ws.Cells.Item(336, 23).Value = 33.33
$ws.Cells.Item(336, 24).Value = 5
$ws.Cells.Item(336, 25).Value = 180.14

# Row 337
$ws.Cells.Item(337, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(337, 2).Value = 8
$ws.Cells.Item(337, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(337, 4).Value = "F. GOMEZ DE ENTERRIA LOPEZ"
$ws.Cells.Item(337, 5).Value = "C.B. TRES CANTOS"
$ws.Cells.Item(337, 6).Value = 1
$ws.Cells.Item(337, 7).Value = 60
$ws.Cells.Item(337, 8).Value = 0
$ws.Cells.Item(337, 9).Value = 1
$ws.Cells.Item(337, 10).Value = 0
$ws.Cells.Item(337, 11).Value = 1
$ws.Cells.Item(337, 12).Value = 0
$ws.Cells.Item(337, 13).Value = 0
$ws.Cells.Item(337, 14).Value = 0
$ws.Cells.Item(337, 15).Value = 0
$ws.Cells.Item(337, 16).Value = 0
$ws.Cells.Item(337, 17).Value = 0
$ws.Cells.Item(337, 18).Value = 0
$ws.Cells.Item(337, 19).Value = 0
$ws.Cells.Item(337, 20).Value = 1
$ws.Cells.Item(337, 21).Value = 1
$ws.Cells.Item(337, 22).Value = 0
$ws.Cells.Item(337, 23).Value = 50
$ws.Cells.Item(337, 24).Value = -2
$ws.Cells.Item(337, 25).Value = -200

# Row 338
$ws.Cells.Item(338, 1).Value = "Liga Regular `"B-B`""
$ws.Cells.Item(338, 2).Value = 8
$ws.Cells.Item(338, 3).Value = "C. D. MENSAJERO ISLA DE LA PALMA vs C.B. TRES CANTOS"
$ws.Cells.Item(338, 4).Value = "A. SANCHO PEREZ"
$ws.Cells.Item(338, 5).Value = "C.B. TRES CANTOS"
$ws.Cells.Item(338, 6).Value = 0.45
$ws.Cells.Item(338, 7).Value = 27
$ws.Cells.Item(338, 8).Value = 0
$ws.Cells.Item(338, 9).Value = 0
$ws.Cells.Item(338, 10).Value = 0
$ws.Cells.Item(338, 11).Value = 0
$ws.Cells.Item(338, 12).Value = 0
$ws.Cells.Item(338, 13).Value = 0
$ws.Cells.Item(338, 14).Value = 0
$ws.Cells.Item(338, 17).Value = 0
$ws.Cells.Item(338, 18).Value = 0
$ws.Cells.Item(338, 19).Value = 0
$ws.Cells.Item(338, 20).Value = 0
$ws.Cells.Item(338, 21).Value = 0
$ws.Cells.Item(338, 22).Value = 0
$ws.Cells.Item(338, 23).Value = 0
$ws.Cells.Item(338, 24).Value = 1
$ws.Cells.Item(338, 25).Value = 72.73

Write-Host "Added rows 326-338"